$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = 10.7957039567556
$ws.Range("B4").Value  = 13212.71166189251
$ws.Range("B5").Value  = 9293.529036356393
$ws.Range("B6").Value  = 920.4767212210036
$ws.Range("B7").Value  = 830.5004175541554
$ws.Range("B8").Value  = 20217.59999999929
$ws.Range("B9").Value  = 4174.859911747636
$ws.Range("B10").Value = 214483.7612266095
$ws.Range("B11").Value = 0.07987641979516172
$ws.Range("B12").Value = 0.450310105782265
$ws.Range("B13").Value = 0.3499999999999948
$ws.Range("B14").Value = 0.9903183921506162
$ws.Range("B15").Value = 0.3943803758639259
